$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-29 Thursday" "2026-01-30 Friday"

Replace-Text "633×4=" "205×5="
Replace-Text "325×9=" "196×5="
Replace-Text "176×7=" "239×9="
Replace-Text "345×8=" "320×8="
Replace-Text "176×6=" "265×6="

Replace-Text "854×8=" "912×8="
Replace-Text "234×5=" "868×8="
Replace-Text "452×9=" "202×5="
Replace-Text "423×7=" "377×3="
Replace-Text "669×4=" "430×6="

Replace-Text "727×6=" "402×3="
Replace-Text "188×4=" "908×9="
Replace-Text "668×6=" "512×6="
Replace-Text "735×3=" "310×9="
Replace-Text "554×2=" "764×5="

Replace-Text "713×2=" "906×7="
Replace-Text "139×2=" "201×7="
Replace-Text "190×9=" "360×9="
Replace-Text "504×4=" "502×4="
Replace-Text "866×6=" "396×2="

Replace-Text "635×9=" "428×4="
Replace-Text "873×9=" "862×3="
Replace-Text "552×7=" "661×9="
Replace-Text "499×6=" "565×6="
Replace-Text "586×7=" "372×8="
